{"js": "const paras = context.document.body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst lastPara = paras.items[paras.items.length - 1];\nconst insertRange = lastPara.getRange(\"After\");\n\nconst newParagraphsXml = `    <w:p/>\n    <w:p>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:t>predictor_variables_recode.dat$</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:t>discipline_father_8y &lt;- recode(predictor_variables_recode.dat$discipline_father_8y, \"0=NA; 2=3; 3=2; 4=2; 5=3; 6=NA\")</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:t>predictor_variables_recode.dat$</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:t>discipline_mother_8y &lt;- recode(predictor_variables_recode.dat$discipline_mother_8y, \"0=NA; 2=3; 3=2; 4=2; 5=3; 6=NA\")</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:t>predictor_variables_recode.dat$</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:t>discipline_quality_father_8y &lt;- recode(predictor_variables_recode.dat$discipline_quality_father_8y, \"0=NA; 4=2; 5=NA\") predictor_variables_recode.dat$discipline_quality_mother_8y &lt;- recode(predictor_variables_recode.dat$discipline_quality_mother_8y, \"0=NA; 4=2; 5=NA\")</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:t>predictor_variables_recode.dat$</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:t>father_interest_8y &lt;- recode(predictor_variables_recode.dat$father_interest_8y, \"0=NA; 3=NA\")</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:t>predictor_variables_recode.dat$</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:t>maternal_attitude_8y &lt;- recode(predictor_variables_recode.dat$maternal_attitude_8y, \"0=NA; 3=2; 4:6=3; 7=NA\")</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:lastRenderedPageBreak/>\n        <w:t>predictor_variables_recode.dat$</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:t>paternal_attitude_8y &lt;- recode(predictor_variables_recode.dat$paternal_attitude_8y, \"0=NA; 4=2; 5:6=NA\")</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:t>predictor_variables_recode.dat$</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:t>neglected_8y &lt;- recode(predictor_variables_recode.dat$neglected_8y, \"0=NA\")</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:t>predictor_variables_recode.dat$</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:t>parent_praise_8y &lt;- recode(predictor_variables_recode.dat$parent_praise_8y, \"0=NA\")</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:t>predictor_variables_recode.dat$</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:t>parent_rules_8y &lt;-recode(predictor_variables_recode.dat$parent_rules_8y, \"0=NA\")</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:t>predictor_variables_recode.dat$</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:t>permanent_sep_parent_8y &lt;- recode(predictor_variables_recode.dat$permanent_sep_parent_8y, \"3=2\")</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:t>predictor_variables_recode.dat$</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:t>over1mo_sep_father_under5y &lt;- recode(predictor_variables_recode.dat$over1mo_sep_father_under5y, \"3=2\")</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:t>predictor_variables_recode.dat$</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:t>over1mo_sep_mother_under5y &lt;- recode(predictor_variables_recode.dat$over1mo_sep_mother_under5y, \"3=2\")</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:t>predictor_variables_recode.dat$</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:t>over1mo_sep_parent_under5y &lt;- recode(predictor_variables_recode.dat$over1mo_sep_parent_under5y, \"3=2\")</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:t>predictor_variables_recode.dat$</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:t>over1mo_sep_father_under10y &lt;- recode(predictor_variables_recode.dat$over1mo_sep_father_under10y, \"3=2\")</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:t>predictor_variables_recode.dat$</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:t>over1mo_sep_mother_under10y &lt;- recode(predictor_variables_recode.dat$over1mo_sep_mother_under10y, \"3=2\")</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:t>predictor_variables_recode.dat$</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:t>over1mo_sep_parent_under10y &lt;- recode(predictor_variables_recode.dat$over1mo_sep_parent_under10y, \"3=2\")</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:t>predictor_variables_recode.dat$</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:t>temp_sep_father_under5y &lt;- recode(predictor_variables_recode.dat$temp_sep_father_under5y, \"3=2\")</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:t>predictor_variables_recode.dat$</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:t>temp_sep_mother_under5y &lt;- recode(predictor_variables_recode.dat$temp_sep_mother_under5y, \"3=2\")</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:t>predictor_variables_recode.dat$</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:t>temp_sep_parent_under5y &lt;-recode(predictor_variables_recode.dat$temp_sep_parent_under5y, \"3=2\")</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:t>predictor_variables_recode.dat$</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:t>temp_sep_father_under10y &lt;- recode(predictor_variables_recode.dat$temp_sep_father_under10y, \"3=2\")</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:t>predictor_variables_recode.dat$</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:t>temp_sep_mother_under10y &lt;- recode(predictor_variables_recode.dat$temp_sep_mother_under10y, \"3=2\")</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:t>predictor_variables_recode.dat$</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:t>temp_sep_parent_under10y &lt;- recode(predictor_variables_recode.dat$temp_sep_parent_under10y, \"3=2\")</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:lastRenderedPageBreak/>\n        <w:t>predictor_variables_recode.dat$</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:t>parent_vigilance_8y &lt;- recode(predictor_variables_recode.dat$parent_vigilance_8y, \"0=NA\")</w:t>\n      </w:r>\n    </w:p>\n    <w:p/>\n    <w:p/>\n    <w:p>\n      <w:r>\n        <w:t>= v159,</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t xml:space=\"preserve\">                             attitude_father_score_8y = v178, attitude_mother_score_8y = v179, pas_authoritarianism_parents_10y = v180, maternal_attitude_10y = v227, paternal_attitude_10y = v248, pas_authoritarianism_father_10y = v249, pas_authoritarianism_mother_10y = v250</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:t>,pas</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:t>_underconcerned_mother_10y = v251, supervision_parents_score_8y = v274,</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t xml:space=\"preserve\">                             parent_approval_12y = v416, permanent_sep_parent_12y = v425,</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t xml:space=\"preserve\">                             live_parents_home_16y = v495,</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t xml:space=\"preserve\">                             parents_alive_18y = v626, live_parentsorfoster_home_16y = v627, want_live_parents_home_18y = v628, agreement_mother_18y = v629, agreement_father_18y = v630, harmony_parents_score_18y = v791,</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:r>\n        <w:t xml:space=\"preserve\">                             age_left_home_21y = v814</w:t>\n\n      </w:r>\n    </w:p>\n`;\n\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n${newParagraphsXml}\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ninsertRange.insertOoxml(ooxml, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$endPos = $d.Content.End\n$rng = $d.Range($endPos, $endPos)\n\n$newParagraphsXml = '<w:p/><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>predictor_variables_recode.dat$</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>discipline_father_8y &lt;- recode(predictor_variables_recode.dat$discipline_father_8y, \"0=NA; 2=3; 3=2; 4=2; 5=3; 6=NA\")</w:t></w:r></w:p><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>predictor_variables_recode.dat$</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>discipline_mother_8y &lt;- recode(predictor_variables_recode.dat$discipline_mother_8y, \"0=NA; 2=3; 3=2; 4=2; 5=3; 6=NA\")</w:t></w:r></w:p><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>predictor_variables_recode.dat$</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>discipline_quality_father_8y &lt;- recode(predictor_variables_recode.dat$discipline_quality_father_8y, \"0=NA; 4=2; 5=NA\") predictor_variables_recode.dat$discipline_quality_mother_8y &lt;- recode(predictor_variables_recode.dat$discipline_quality_mother_8y, \"0=NA; 4=2; 5=NA\")</w:t></w:r></w:p><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>predictor_variables_recode.dat$</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>father_interest_8y &lt;- recode(predictor_variables_recode.dat$father_interest_8y, \"0=NA; 3=NA\")</w:t></w:r></w:p><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>predictor_variables_recode.dat$</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>maternal_attitude_8y &lt;- recode(predictor_variables_recode.dat$maternal_attitude_8y, \"0=NA; 3=2; 4:6=3; 7=NA\")</w:t></w:r></w:p><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:lastRenderedPageBreak/><w:t>predictor_variables_recode.dat$</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>paternal_attitude_8y &lt;- recode(predictor_variables_recode.dat$paternal_attitude_8y, \"0=NA; 4=2; 5:6=NA\")</w:t></w:r></w:p><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>predictor_variables_recode.dat$</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>neglected_8y &lt;- recode(predictor_variables_recode.dat$neglected_8y, \"0=NA\")</w:t></w:r></w:p><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>predictor_variables_recode.dat$</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>parent_praise_8y &lt;- recode(predictor_variables_recode.dat$parent_praise_8y, \"0=NA\")</w:t></w:r></w:p><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>predictor_variables_recode.dat$</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>parent_rules_8y &lt;-recode(predictor_variables_recode.dat$parent_rules_8y, \"0=NA\")</w:t></w:r></w:p><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>predictor_variables_recode.dat$</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>permanent_sep_parent_8y &lt;- recode(predictor_variables_recode.dat$permanent_sep_parent_8y, \"3=2\")</w:t></w:r></w:p><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>predictor_variables_recode.dat$</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>over1mo_sep_father_under5y &lt;- recode(predictor_variables_recode.dat$over1mo_sep_father_under5y, \"3=2\")</w:t></w:r></w:p><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>predictor_variables_recode.dat$</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>over1mo_sep_mother_under5y &lt;- recode(predictor_variables_recode.dat$over1mo_sep_mother_under5y, \"3=2\")</w:t></w:r></w:p><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>predictor_variables_recode.dat$</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>over1mo_sep_parent_under5y &lt;- recode(predictor_variables_recode.dat$over1mo_sep_parent_under5y, \"3=2\")</w:t></w:r></w:p><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>predictor_variables_recode.dat$</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>over1mo_sep_father_under10y &lt;- recode(predictor_variables_recode.dat$over1mo_sep_father_under10y, \"3=2\")</w:t></w:r></w:p><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>predictor_variables_recode.dat$</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>over1mo_sep_mother_under10y &lt;- recode(predictor_variables_recode.dat$over1mo_sep_mother_under10y, \"3=2\")</w:t></w:r></w:p><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>predictor_variables_recode.dat$</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>over1mo_sep_parent_under10y &lt;- recode(predictor_variables_recode.dat$over1mo_sep_parent_under10y, \"3=2\")</w:t></w:r></w:p><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>predictor_variables_recode.dat$</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>temp_sep_father_under5y &lt;- recode(predictor_variables_recode.dat$temp_sep_father_under5y, \"3=2\")</w:t></w:r></w:p><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>predictor_variables_recode.dat$</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>temp_sep_mother_under5y &lt;- recode(predictor_variables_recode.dat$temp_sep_mother_under5y, \"3=2\")</w:t></w:r></w:p><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>predictor_variables_recode.dat$</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>temp_sep_parent_under5y &lt;-recode(predictor_variables_recode.dat$temp_sep_parent_under5y, \"3=2\")</w:t></w:r></w:p><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>predictor_variables_recode.dat$</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>temp_sep_father_under10y &lt;- recode(predictor_variables_recode.dat$temp_sep_father_under10y, \"3=2\")</w:t></w:r></w:p><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>predictor_variables_recode.dat$</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>temp_sep_mother_under10y &lt;- recode(predictor_variables_recode.dat$temp_sep_mother_under10y, \"3=2\")</w:t></w:r></w:p><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>predictor_variables_recode.dat$</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>temp_sep_parent_under10y &lt;- recode(predictor_variables_recode.dat$temp_sep_parent_under10y, \"3=2\")</w:t></w:r></w:p><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:lastRenderedPageBreak/><w:t>predictor_variables_recode.dat$</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>parent_vigilance_8y &lt;- recode(predictor_variables_recode.dat$parent_vigilance_8y, \"0=NA\")</w:t></w:r></w:p><w:p/><w:p/><w:p><w:r><w:t>= v159,</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">                             attitude_father_score_8y = v178, attitude_mother_score_8y = v179, pas_authoritarianism_parents_10y = v180, maternal_attitude_10y = v227, paternal_attitude_10y = v248, pas_authoritarianism_father_10y = v249, pas_authoritarianism_mother_10y = v250</w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:t>,pas</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>_underconcerned_mother_10y = v251, supervision_parents_score_8y = v274,</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">                             parent_approval_12y = v416, permanent_sep_parent_12y = v425,</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">                             live_parents_home_16y = v495,</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">                             parents_alive_18y = v626, live_parentsorfoster_home_16y = v627, want_live_parents_home_18y = v628, agreement_mother_18y = v629, agreement_father_18y = v630, harmony_parents_score_18y = v791,</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">                             age_left_home_21y = v814</w:t></w:r></w:p>'\n\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + $newParagraphsXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$rng.InsertXML($xml)\n"}
